$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5823119999999999
$ws.Range("H2").Value = 1.746936
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.21127
$ws.Range("N2").Value = 0.63381
$ws.Range("O2").Value = 0.02468048274231428
$ws.Range("P2").Value = 0.02468048274231428
$ws.Range("Q2").Value = 0.12302505624
$ws.Range("R2").Value = 1.10722550616
$ws.Range("S2").Value = 0.02468048274231428
$ws.Range("T2").Value = 0.02468048274231428

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5823119999999999
$ws.Range("H3").Value = 1.746936
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.467027333333334
$ws.Range("N3").Value = 10.401082
$ws.Range("O3").Value = 0.405016842275123
$ws.Range("P3").Value = 0.405016842275123
$ws.Range("Q3").Value = 2.018891620528
$ws.Range("R3").Value = 18.170024584752
$ws.Range("S3").Value = 0.405016842275123
$ws.Range("T3").Value = 0.405016842275123

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.5823119999999999
$ws.Range("H4").Value = 1.746936
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.881908
$ws.Range("N4").Value = 14.645724
$ws.Range("O4").Value = 0.5703026749825627
$ws.Range("P4").Value = 0.5703026749825627
$ws.Range("Q4").Value = 2.842793611296
$ws.Range("R4").Value = 25.585142501664
$ws.Range("S4").Value = 0.5703026749825627
$ws.Range("T4").Value = 0.5703026749825627
